$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C (y_0_forecast) and Column E (y_1_forecast) updated values per row
$values = @{
    2  = @{ C = 5.021907707863549;   E = 19.19812743658083 }
    3  = @{ C = -5.592633745595466;  E = -9.964084247724703 }
    4  = @{ C = -5.440152375872254;  E = -14.43639438706738 }
    5  = @{ C = 9.349082908138451;   E = 27.15801420548429 }
    6  = @{ C = 0.5389546843750148;  E = -5.27893918837793 }
    7  = @{ C = -4.232836797447693;  E = -8.0930759205322 }
    8  = @{ C = 7.942828065321739;   E = 15.37760125310905 }
    9  = @{ C = 1.913895196850168;   E = 6.97490799213798 }
    10 = @{ C = 4.861901970953975;   E = 9.631040506010535 }
    11 = @{ C = 4.115488239647713;   E = 9.523050046161053 }
    12 = @{ C = 4.07381142256642;    E = 8.54956688663686 }
    13 = @{ C = 5.264109583376908;   E = 7.819356632099961 }
    14 = @{ C = 6.942957493752444;   E = 17.43645097609996 }
    15 = @{ C = 2.932994663878907;   E = 2.894715150804616 }
    16 = @{ C = -0.7538332529782865; E = -4.784481399264983 }
    17 = @{ C = -2.598185084325777;  E = -3.071148328823314 }
    18 = @{ C = -1.287508943286542;  E = -4.200823682253607 }
    19 = @{ C = 0.1893861904177951;  E = 0.1126281723122791 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
